$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "D-SUB 15 pol HD hona" hyperlink (on D4) is being retired - the row's
# text/url get a "PCB" variant below, but it no longer carries a live
# hyperlink relationship. Remove just that one hyperlink (leaves D2/D3
# hyperlinks and their formatting untouched).
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Address -eq "https://www.electrokit.com/d-sub-15-pol-hd-hona") {
        $hl.Delete()
    }
}

# New row 5: LCD TFT 2.8" touch (ILI9341), amount 1
$ws.Range("D5").Value = "https://www.electrokit.com/lcd-tft-2.8touch-ili9341"
$ws.Range("A5").Value = "LCD TFT 2.8"" touch (ILI9341)"
$ws.Range("B5").Value = 1

# New row 6: 6.3mm hona 2-pol mono chassie brytare, amount 2
$ws.Range("D6").Value = "https://www.electrokit.com/6.3mm-chassie-mono-med-brytare"
$ws.Range("A6").Value = "6.3mm hona 2-pol mono chassie brytare"
$ws.Range("B6").Value = 2

# Row 4: "D-SUB 15 pol HD hona" -> "D-SUB 15 pol HD hona PCB" text/url update
$ws.Range("D4").Value = "https://www.electrokit.com/d-sub-15-pol-hd-hona-pcb"
$ws.Range("A4").Value = "D-SUB 15 pol HD hona PCB"

# Wire up hyperlinks for the new rows
$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.electrokit.com/lcd-tft-2.8touch-ili9341")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://www.electrokit.com/6.3mm-chassie-mono-med-brytare")

# Hyperlinks.Add leaves behind slightly different direct formatting; re-apply
# the workbook's "Hyperlänk" cell style so D4/D5/D6 match D2/D3 exactly.
$ws.Range("D4").ClearFormats()
$ws.Range("D4").Style = "Hyperlänk"
$ws.Range("D5").ClearFormats()
$ws.Range("D5").Style = "Hyperlänk"
$ws.Range("D6").ClearFormats()
$ws.Range("D6").Style = "Hyperlänk"

Write-Host "done"
